$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (A2:D2 = 0, 0, 0.01, 0) which shifts rows 3:8 up to 2:7.
$ws.Rows.Item(2).Delete()

# Update the selection/active cell to D8, matching the post-edit state.
$ws.Range("D8").Select()

$wb.Save()
